$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 and Row 15 swap values in columns A, Q, R
$a14 = $ws.Range("A14").Value2
$q14 = $ws.Range("Q14").Value2
$r14 = $ws.Range("R14").Value2

$a15 = $ws.Range("A15").Value2
$q15 = $ws.Range("Q15").Value2
$r15 = $ws.Range("R15").Value2

$ws.Range("A14").Value2 = $a15
$ws.Range("Q14").Value2 = $q15
$ws.Range("R14").Value2 = $r15

$ws.Range("A15").Value2 = $a14
$ws.Range("Q15").Value2 = $q14
$ws.Range("R15").Value2 = $r14

# Row 16 and Row 17 swap values in columns A, Q, R, AC
$a16 = $ws.Range("A16").Value2
$q16 = $ws.Range("Q16").Value2
$r16 = $ws.Range("R16").Value2
$ac16 = $ws.Range("AC16").Value2

$a17 = $ws.Range("A17").Value2
$q17 = $ws.Range("Q17").Value2
$r17 = $ws.Range("R17").Value2
$ac17 = $ws.Range("AC17").Value2

$ws.Range("A16").Value2 = $a17
$ws.Range("Q16").Value2 = $q17
$ws.Range("R16").Value2 = $r17
$ws.Range("AC16").Value2 = $ac17

$ws.Range("A17").Value2 = $a16
$ws.Range("Q17").Value2 = $q16
$ws.Range("R17").Value2 = $r16
$ws.Range("AC17").Value2 = $ac16
